# Append three new invoice/posting rows (61-63) at the bottom of the
# "bd_factures_fournisseurs" data table, pushing the sheet's dimension
# from A1:X60 to A1:X63.
#
# Every value in the source data is plain text (dates, amounts and
# sequence numbers are all stored as literal strings, e.g. "2025-07-24",
# "132.00", "23"), so each cell is written with NumberFormat "@" first to
# stop Excel from auto-converting the text into a real date/number, and
# then ClearFormats() removes the now-useless "@" format again so the
# cell keeps the workbook's default (unstyled) look - only the literal
# text sticks around. Cells that must stay an explicit empty string
# (as opposed to a truly blank/absent cell) are written with a lone
# apostrophe formula (forces a Text-typed empty cell) and then also
# cleared of the quote-prefix formatting that produces.

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

function Set-EmptyTextCell($range) {
    $range.Formula = "'"
    $range.ClearFormats()
}

function Set-RowValues($ws, $rowNumber, $values) {
    foreach ($col in $values.Keys) {
        $address = "$col$rowNumber"
        $range = $ws.Range($address)
        $value = $values[$col]
        if ($value -eq $null) {
            Set-EmptyTextCell $range
        } else {
            Set-TextCell $range $value
        }
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        Row    = 61
        Values = [ordered]@{
            A = "Fournisseur_02"
            B = "1001 – Caisse 2"
            C = "31 jours"
            D = "2025-07-24"
            E = "2025-07-24"
            F = "2025-07-24"
            G = "2025-07"
            H = "23"
            I = "34"
            J = "1"
            K = "2"
            L = "3"
            M = "1010 – CCP 1"
            N = "test24"
            O = "44"
            P = "3"
            Q = "132.00"
            R = "1011 – CCP 2"
            S = "test24"
            T = "7.7"
            U = "4"
            V = "136.00"
            W = $null
            X = "1012"
        }
    },
    @{
        Row    = 62
        Values = [ordered]@{
            A = "Fournisseur_07"
            B = "1100 – Débiteurs 1"
            C = "36 jours"
            D = "2025-07-25"
            E = $null
            F = $null
            G = "2025-07"
            H = "3"
            I = "4"
            J = "5"
            K = "6"
            L = "7"
            M = "1101 – Débiteurs 2"
            N = "8"
            O = "9"
            P = "101"
            Q = "909.00"
            R = "1200 – Stock 1"
            S = "11"
            T = "7.7"
            U = "12"
            V = "921.00"
            W = $null
            X = "1013"
        }
    },
    @{
        Row    = 63
        Values = [ordered]@{
            A = "Fournisseur_03"
            B = "1010 – CCP 1"
            C = "32 jours"
            D = "2025-07-25"
            E = $null
            F = $null
            G = "2025-07"
            H = "26"
            I = "27"
            J = "28"
            K = "29"
            L = "30"
            M = "1011 – CCP 2"
            N = "31"
            O = "32"
            P = "33"
            Q = "1056.00"
            R = "1020 – Banque 1"
            S = "34"
            T = "7.7"
            U = "35"
            V = "1091.00"
            W = $null
            X = "1014"
        }
    }
)

foreach ($rowDef in $newRows) {
    Set-RowValues $ws $rowDef.Row $rowDef.Values
}
